# Add three new joiner rows (109-111) to the profiles sheet, matching the
# name/username/flag layout used throughout the table, then update the
# view's scroll position and selection to land on the newly added data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 109: Prabhat Kumar / prabhat236218 -------------------------------
# Values are written in the same left-to-right order the source file shows
# in the shared-string table (name, then username, then the flag).
$ws.Range("A109").Value = "Prabhat Kumar"
$ws.Range("B109").Value = "prabhat236218"
$ws.Range("C109").Value = 0

# --- Row 110: Raghu Mahajan / raghu96 -------------------------------------
# Username (column B) is entered before the name (column A) here, matching
# the original authoring order.
$ws.Range("B110").Value = "raghu96"
$ws.Range("A110").Value = "Raghu Mahajan"
$ws.Range("C110").Value = 0

# --- Row 111: Hardy Tom / hardytom ----------------------------------------
$ws.Range("A111").Value = "Hardy Tom"
$ws.Range("B111").Value = "hardytom"
$ws.Range("C111").Value = 0

# Match the formatting used by the closest existing rows with the same
# per-column style pattern (A populated, B plain, C flagged) so the new
# entries look consistent with the rest of the table.
$ws.Range("A99:C99").Copy()
$ws.Range("A109:C109").PasteSpecial(-4122)

$ws.Range("A106:C106").Copy()
$ws.Range("A110:C110").PasteSpecial(-4122)

$ws.Range("A107:C107").Copy()
$ws.Range("A111:C111").PasteSpecial(-4122)

# Scroll the sheet so row 89 is at the top and select the last entered cell,
# mirroring where the editor was working.
$excel.ActiveWindow.ScrollRow = 89
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C111").Select()
